$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.239.44'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '3.352.96'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '583.25'
$c.ClearFormats()
$ws.Range("E5").Value = '  +0.27%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '177.05'
$c.ClearFormats()
$ws.Range("E6").Value = '  +0.37%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range("E7").Value = '  +0.06%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.590'
$c.ClearFormats()
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +2.09%  '
$ws.Range("E10").Value = '  +0.69%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '48.06'
$c.ClearFormats()
$ws.Range("E11").Value = '  +5.53%  '
$ws.Range("E12").Value = '  +1.41%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '685.33'
$c.ClearFormats()
$ws.Range("E13").Value = '  +3.88%  '
$ws.Range("D14").Value = '3.890.19'
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").Value = '68.255.01'
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").Value = '3.349.39'
$ws.Range("E18").Value = '  +0.71%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '17.45'
$c.ClearFormats()
$ws.Range("E19").Value = '  +0.34%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '11.20'
$c.ClearFormats()
$ws.Range("E20").Value = '  +2.12%  '
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("E22").Value = '  -0.47%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '16.97'
$c.ClearFormats()
$ws.Range("E23").Value = '  -0.70%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '100.42'
$c.ClearFormats()
$ws.Range("E24").Value = '  +1.10%  '
$ws.Range("E25").Value = '  +1.40%  '
$ws.Range("E26").Value = '  +0.85%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.51'
$c.ClearFormats()
$ws.Range("E27").Value = '  +2.39%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '33.00'
$c.ClearFormats()
$ws.Range("E28").Value = '  -1.87%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '8.52'
$c.ClearFormats()
$ws.Range("E29").Value = '  +0.81%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '6.94'
$c.ClearFormats()
$ws.Range("E30").Value = '  -7.16%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '561.57'
$c.ClearFormats()
$ws.Range("E31").Value = '  -2.51%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '11.08'
$c.ClearFormats()
$ws.Range("E32").Value = '  +0.93%  '
$ws.Range("E33").Value = '  +1.07%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '57.98'
$c.ClearFormats()
$ws.Range("E34").Value = '  +2.67%  '
$ws.Range("D36").Value = '3.719.52'
$ws.Range("E36").Value = '  +0.61%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.32'
$c.ClearFormats()
$ws.Range("E37").Value = '  -2.32%  '
$ws.Range("E38").Value = '  +4.92%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '34.82'
$c.ClearFormats()
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("E40").Value = '  +1.52%  '
$ws.Range("E41").Value = '  -0.43%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.336'
$c.ClearFormats()
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("D43").Value = '0.0₃0672'
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("E44").Value = '  -1.23%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0411'
$c.ClearFormats()
$ws.Range("E45").Value = '  +1.21%  '
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("E49").Value = '  -0.49%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '131.88'
$c.ClearFormats()
$ws.Range("E50").Value = '  +3.07%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.57'
$c.ClearFormats()
$ws.Range("E51").Value = '  -1.77%  '
